# Adds a newer bank-statement movement ("INTERES A SU FAVOR" dated
# 2014-02-03) above the existing rows by inserting a new row 1, which
# pushes all prior rows (and the CONCATENATE() formula that lived on the
# old row 1) down by one, and then fills in the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 1; everything below shifts
# down (row 1 -> row 2, row 2 -> row 3, ... and the trailing blank rows
# grow by one too, e.g. old row 28 -> new row 29).
$ws.Rows.Item(1).Insert()

# The inserted row starts out with "General" formatting on column A, but
# every other row in this sheet shows the date as m/d/yyyy. Pull that
# formatting down from the row right below (the row that used to be row 1)
# so the new A1 renders as a date instead of a bare serial number.
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Non-breaking space used inside the padded amount strings in this sheet
# (e.g. "0.25  " actually contains U+00A0 twice, not regular spaces).
$nbsp = [char]0x00A0

$ws.Range("A1").Value = 41673
$ws.Range("B1").Value = "INTERES A SU FAVOR"
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "0000952078"
$ws.Range("E1").Value = "AGENCIA PARA PROCESOS BATCH"
$ws.Range("F1").Value = "0.25$nbsp$nbsp"
$ws.Range("G1").Value = "3992.30"
